# Add a new "Total Taxes (% GDP)" column (M) to the cleaned totals table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header (row 1) ---
$ws.Range("M1").Value = "Total Taxes (% GDP)"
$ws.Range("M1").NumberFormat = "0.00"
$ws.Range("M1").Font.Bold = $true

# --- Data rows (2-148) ---
$ws.Range("M2").Value = 9.1248640000000005
$ws.Range("M2").NumberFormat = "0.00"
$ws.Range("M3").Value = 23.52533
$ws.Range("M3").NumberFormat = "0.00"
$ws.Range("M4").Value = 14.162544
$ws.Range("M4").NumberFormat = "0.00"
$ws.Range("M5").Value = 30.661214999999999
$ws.Range("M5").NumberFormat = "0.00"
$ws.Range("M6").Value = 21.761512
$ws.Range("M6").NumberFormat = "0.00"
$ws.Range("M7").Value = 28.578520000000001
$ws.Range("M7").NumberFormat = "0.00"
$ws.Range("M8").Value = 41.884256999999998
$ws.Range("M8").NumberFormat = "0.00"
$ws.Range("M9").Value = 16.506067999999999
$ws.Range("M9").NumberFormat = "0.00"
$ws.Range("M10").Value = 1.0698810000000001
$ws.Range("M10").NumberFormat = "0.00"
$ws.Range("M11").Value = 8.1908390000000004
$ws.Range("M11").NumberFormat = "0.00"
$ws.Range("M12").Value = 35.748614000000003
$ws.Range("M12").NumberFormat = "0.00"
$ws.Range("M13").Value = 43.875816
$ws.Range("M13").NumberFormat = "0.00"
$ws.Range("M14").Value = 15.888945
$ws.Range("M14").NumberFormat = "0.00"
$ws.Range("M15").Value = 26.203921000000001
$ws.Range("M15").NumberFormat = "0.00"
$ws.Range("M16").Value = 37.572918999999999
$ws.Range("M16").NumberFormat = "0.00"
$ws.Range("M17").Value = 20.871931
$ws.Range("M17").NumberFormat = "0.00"
$ws.Range("M18").Value = 33.713017999999998
$ws.Range("M18").NumberFormat = "0.00"
$ws.Range("M19").Value = 28.800666
$ws.Range("M19").NumberFormat = "0.00"
$ws.Range("M20").Value = 16.652574000000001
$ws.Range("M20").NumberFormat = "0.00"
$ws.Range("M21").Value = 14.825889999999999
$ws.Range("M21").NumberFormat = "0.00"
$ws.Range("M22").Value = 13.376192
$ws.Range("M22").NumberFormat = "0.00"
$ws.Range("M23").Value = 32.830274000000003
$ws.Range("M23").NumberFormat = "0.00"
$ws.Range("M24").Value = 9.6534659999999999
$ws.Range("M24").NumberFormat = "0.00"
$ws.Range("M25").Value = 5.6650309999999999
$ws.Range("M25").NumberFormat = "0.00"
$ws.Range("M26").Value = 17.409112
$ws.Range("M26").NumberFormat = "0.00"
$ws.Range("M27").Value = 23.790354000000001
$ws.Range("M27").NumberFormat = "0.00"
$ws.Range("M28").Value = 19.211746000000002
$ws.Range("M28").NumberFormat = "0.00"
$ws.Range("M29").Value = 17.503682999999999
$ws.Range("M29").NumberFormat = "0.00"
$ws.Range("M30").Value = 7.1134009999999996
$ws.Range("M30").NumberFormat = "0.00"
$ws.Range("M31").Value = 13.854123
$ws.Range("M31").NumberFormat = "0.00"
$ws.Range("M32").Value = 35.163044999999997
$ws.Range("M32").NumberFormat = "0.00"
$ws.Range("M33").Value = 23.997837000000001
$ws.Range("M33").NumberFormat = "0.00"
$ws.Range("M34").Value = 34.068710000000003
$ws.Range("M34").NumberFormat = "0.00"
$ws.Range("M35").Value = 45.43515
$ws.Range("M35").NumberFormat = "0.00"
$ws.Range("M36").Value = 13.509418
$ws.Range("M36").NumberFormat = "0.00"
$ws.Range("M37").Value = 14.671663000000001
$ws.Range("M37").NumberFormat = "0.00"
$ws.Range("M38").NumberFormat = "0.00"
$ws.Range("M39").Value = 17.381668999999999
$ws.Range("M39").NumberFormat = "0.00"
$ws.Range("M40").Value = 33.730350999999999
$ws.Range("M40").NumberFormat = "0.00"
$ws.Range("M41").NumberFormat = "0.00"
$ws.Range("M42").Value = 44.017195000000001
$ws.Range("M42").NumberFormat = "0.00"
$ws.Range("M43").Value = 45.277869000000003
$ws.Range("M43").NumberFormat = "0.00"
$ws.Range("M44").Value = 20.494221
$ws.Range("M44").NumberFormat = "0.00"
$ws.Range("M45").Value = 11.14148
$ws.Range("M45").NumberFormat = "0.00"
$ws.Range("M46").Value = 25.819448000000001
$ws.Range("M46").NumberFormat = "0.00"
$ws.Range("M47").Value = 37.430619999999998
$ws.Range("M47").NumberFormat = "0.00"
$ws.Range("M48").Value = 12.377316
$ws.Range("M48").NumberFormat = "0.00"
$ws.Range("M49").Value = 38.338583999999997
$ws.Range("M49").NumberFormat = "0.00"
$ws.Range("M50").Value = 10.78232
$ws.Range("M50").NumberFormat = "0.00"
$ws.Range("M51").Value = 14.030989
$ws.Range("M51").NumberFormat = "0.00"
$ws.Range("M52").Value = 13.731528000000001
$ws.Range("M52").NumberFormat = "0.00"
$ws.Range("M53").Value = 22.421873000000001
$ws.Range("M53").NumberFormat = "0.00"
$ws.Range("M54").Value = 13.574933
$ws.Range("M54").NumberFormat = "0.00"
$ws.Range("M55").Value = 39.148890999999999
$ws.Range("M55").NumberFormat = "0.00"
$ws.Range("M56").Value = 36.789971000000001
$ws.Range("M56").NumberFormat = "0.00"
$ws.Range("M57").Value = 17.577255000000001
$ws.Range("M57").NumberFormat = "0.00"
$ws.Range("M58").Value = 10.361459
$ws.Range("M58").NumberFormat = "0.00"
$ws.Range("M59").Value = 7.9746290000000002
$ws.Range("M59").NumberFormat = "0.00"
$ws.Range("M60").Value = 1.3778090000000001
$ws.Range("M60").NumberFormat = "0.00"
$ws.Range("M61").Value = 23.358557999999999
$ws.Range("M61").NumberFormat = "0.00"
$ws.Range("M62").Value = 31.128926
$ws.Range("M62").NumberFormat = "0.00"
$ws.Range("M63").Value = 42.385586000000004
$ws.Range("M63").NumberFormat = "0.00"
$ws.Range("M64").Value = 16.017084000000001
$ws.Range("M64").NumberFormat = "0.00"
$ws.Range("M65").Value = 25.114878999999998
$ws.Range("M65").NumberFormat = "0.00"
$ws.Range("M66").Value = 30.777436000000002
$ws.Range("M66").NumberFormat = "0.00"
$ws.Range("M67").Value = 15.341718
$ws.Range("M67").NumberFormat = "0.00"
$ws.Range("M68").Value = 15.228497000000001
$ws.Range("M68").NumberFormat = "0.00"
$ws.Range("M69").Value = 15.800242000000001
$ws.Range("M69").NumberFormat = "0.00"
$ws.Range("M70").Value = 23.459638000000002
$ws.Range("M70").NumberFormat = "0.00"
$ws.Range("M71").Value = 1.5125850000000001
$ws.Range("M71").NumberFormat = "0.00"
$ws.Range("M72").Value = 25.102250000000002
$ws.Range("M72").NumberFormat = "0.00"
$ws.Range("M73").Value = 12.428921000000001
$ws.Range("M73").NumberFormat = "0.00"
$ws.Range("M74").Value = 30.234611000000001
$ws.Range("M74").NumberFormat = "0.00"
$ws.Range("M75").Value = 13.693778999999999
$ws.Range("M75").NumberFormat = "0.00"
$ws.Range("M76").Value = 26.536635
$ws.Range("M76").NumberFormat = "0.00"
$ws.Range("M77").NumberFormat = "0.00"
$ws.Range("M78").Value = 1.192523
$ws.Range("M78").NumberFormat = "0.00"
$ws.Range("M79").Value = 29.684875000000002
$ws.Range("M79").NumberFormat = "0.00"
$ws.Range("M80").Value = 37.932749000000001
$ws.Range("M80").NumberFormat = "0.00"
$ws.Range("M81").Value = 25.286065000000001
$ws.Range("M81").NumberFormat = "0.00"
$ws.Range("M82").Value = 10.999344000000001
$ws.Range("M82").NumberFormat = "0.00"
$ws.Range("M83").Value = 15.323005
$ws.Range("M83").NumberFormat = "0.00"
$ws.Range("M84").Value = 14.916247
$ws.Range("M84").NumberFormat = "0.00"
$ws.Range("M85").Value = 31.893222000000002
$ws.Range("M85").NumberFormat = "0.00"
$ws.Range("M86").Value = 16.519642999999999
$ws.Range("M86").NumberFormat = "0.00"
$ws.Range("M87").Value = 18.213287999999999
$ws.Range("M87").NumberFormat = "0.00"
$ws.Range("M88").Value = 13.701247
$ws.Range("M88").NumberFormat = "0.00"
$ws.Range("M89").Value = 26.432732999999999
$ws.Range("M89").NumberFormat = "0.00"
$ws.Range("M90").Value = 20.360968
$ws.Range("M90").NumberFormat = "0.00"
$ws.Range("M91").Value = 36.090035
$ws.Range("M91").NumberFormat = "0.00"
$ws.Range("M92").Value = 21.448879999999999
$ws.Range("M92").NumberFormat = "0.00"
$ws.Range("M93").Value = 20.156713
$ws.Range("M93").NumberFormat = "0.00"
$ws.Range("M94").Value = 8.3322769999999995
$ws.Range("M94").NumberFormat = "0.00"
$ws.Range("M95").Value = 28.710813000000002
$ws.Range("M95").NumberFormat = "0.00"
$ws.Range("M96").Value = 18.687477999999999
$ws.Range("M96").NumberFormat = "0.00"
$ws.Range("M97").Value = 38.400931
$ws.Range("M97").NumberFormat = "0.00"
$ws.Range("M98").Value = 32.613869999999999
$ws.Range("M98").NumberFormat = "0.00"
$ws.Range("M99").Value = 23.448450000000001
$ws.Range("M99").NumberFormat = "0.00"
$ws.Range("M100").Value = 13.454078000000001
$ws.Range("M100").NumberFormat = "0.00"
$ws.Range("M101").Value = 7.1707000000000001
$ws.Range("M101").NumberFormat = "0.00"
$ws.Range("M102").Value = 38.709806999999998
$ws.Range("M102").NumberFormat = "0.00"
$ws.Range("M103").Value = 9.9883570000000006
$ws.Range("M103").NumberFormat = "0.00"
$ws.Range("M104").Value = 5.5946610000000003
$ws.Range("M104").NumberFormat = "0.00"
$ws.Range("M105").Value = 15.463361000000001
$ws.Range("M105").NumberFormat = "0.00"
$ws.Range("M106").Value = 10.489746
$ws.Range("M106").NumberFormat = "0.00"
$ws.Range("M107").Value = 13.614037
$ws.Range("M107").NumberFormat = "0.00"
$ws.Range("M108").Value = 13.676774
$ws.Range("M108").NumberFormat = "0.00"
$ws.Range("M109").Value = 33.346711999999997
$ws.Range("M109").NumberFormat = "0.00"
$ws.Range("M110").Value = 34.141741000000003
$ws.Range("M110").NumberFormat = "0.00"
$ws.Range("M111").Value = 26.441125
$ws.Range("M111").NumberFormat = "0.00"
$ws.Range("M112").Value = 29.046813
$ws.Range("M112").NumberFormat = "0.00"
$ws.Range("M113").Value = 15.459153000000001
$ws.Range("M113").NumberFormat = "0.00"
$ws.Range("M114").Value = 2.23278
$ws.Range("M114").NumberFormat = "0.00"
$ws.Range("M115").Value = 15.920230999999999
$ws.Range("M115").NumberFormat = "0.00"
$ws.Range("M116").Value = 35.078639000000003
$ws.Range("M116").NumberFormat = "0.00"
$ws.Range("M117").Value = 11.178661
$ws.Range("M117").NumberFormat = "0.00"
$ws.Range("M118").Value = 13.421915
$ws.Range("M118").NumberFormat = "0.00"
$ws.Range("M119").Value = 32.356543000000002
$ws.Range("M119").NumberFormat = "0.00"
$ws.Range("M120").Value = 36.505336999999997
$ws.Range("M120").NumberFormat = "0.00"
$ws.Range("M121").Value = 29.305698
$ws.Range("M121").NumberFormat = "0.00"
$ws.Range("M122").Value = 26.236792000000001
$ws.Range("M122").NumberFormat = "0.00"
$ws.Range("M123").NumberFormat = "0.00"
$ws.Range("M124").Value = 33.180988999999997
$ws.Range("M124").NumberFormat = "0.00"
$ws.Range("M125").Value = 12.444495999999999
$ws.Range("M125").NumberFormat = "0.00"
$ws.Range("M126").Value = 44.237478000000003
$ws.Range("M126").NumberFormat = "0.00"
$ws.Range("M127").Value = 27.700901999999999
$ws.Range("M127").NumberFormat = "0.00"
$ws.Range("M128").NumberFormat = "0.00"
$ws.Range("M129").Value = 20.6
$ws.Range("M129").NumberFormat = "0.00"
$ws.Range("M130").Value = 11.390269
$ws.Range("M130").NumberFormat = "0.00"
$ws.Range("M131").Value = 17.834917999999998
$ws.Range("M131").NumberFormat = "0.00"
$ws.Range("M132").Value = 18.763441
$ws.Range("M132").NumberFormat = "0.00"
$ws.Range("M133").Value = 26.147864999999999
$ws.Range("M133").NumberFormat = "0.00"
$ws.Range("M134").Value = 29.614932
$ws.Range("M134").NumberFormat = "0.00"
$ws.Range("M135").Value = 25.301887000000001
$ws.Range("M135").NumberFormat = "0.00"
$ws.Range("M136").Value = 20.373964000000001
$ws.Range("M136").NumberFormat = "0.00"
$ws.Range("M137").Value = 12.676283
$ws.Range("M137").NumberFormat = "0.00"
$ws.Range("M138").Value = 30.859074
$ws.Range("M138").NumberFormat = "0.00"
$ws.Range("M139").Value = 9.2430889999999994
$ws.Range("M139").NumberFormat = "0.00"
$ws.Range("M140").Value = 32.734135000000002
$ws.Range("M140").NumberFormat = "0.00"
$ws.Range("M141").Value = 25.773724000000001
$ws.Range("M141").NumberFormat = "0.00"
$ws.Range("M142").Value = 34.229075999999999
$ws.Range("M142").NumberFormat = "0.00"
$ws.Range("M143").Value = 26.204868999999999
$ws.Range("M143").NumberFormat = "0.00"
$ws.Range("M144").Value = 20.279105000000001
$ws.Range("M144").NumberFormat = "0.00"
$ws.Range("M145").Value = 17.944656999999999
$ws.Range("M145").NumberFormat = "0.00"
$ws.Range("M146").Value = 7.0387729999999999
$ws.Range("M146").NumberFormat = "0.00"
$ws.Range("M147").Value = 12.970041
$ws.Range("M147").NumberFormat = "0.00"
$ws.Range("M148").Value = 27.182075999999999
$ws.Range("M148").NumberFormat = "0.00"

# Restore the selection/scroll position shown in the saved workbook.
$ws.Range("F144").Select()
